$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -0.3538576297953431
$ws.Cells.Item(2, 3).Value = 0.4918198143377789
$ws.Cells.Item(2, 4).Value = 0.02616020827444671
$ws.Cells.Item(2, 5).Value = 0.7056676249657313
$ws.Cells.Item(2, 6).Value = 0.9099227904698376
$ws.Cells.Item(2, 7).Value = 0.4696485924218674
$ws.Cells.Item(2, 8).Value = 0.3872693173527261
$ws.Cells.Item(2, 9).Value = 0.6785207862006621
$ws.Cells.Item(2, 10).Value = 0.5031320576645155
$ws.Cells.Item(2, 11).Value = 0.6535938957934404
$ws.Cells.Item(3, 2).Value = 0.5142477939707875
$ws.Cells.Item(3, 3).Value = 0.02876374783689539
$ws.Cells.Item(3, 4).Value = 0.683599939311716
$ws.Cells.Item(3, 5).Value = 0.9039408101654827
$ws.Cells.Item(3, 6).Value = 0.4643053422394562
$ws.Cells.Item(3, 7).Value = 0.3769007839823598
$ws.Cells.Item(3, 8).Value = 0.6700434640770384
$ws.Cells.Item(3, 9).Value = 0.4953597493577047
$ws.Cells.Item(3, 10).Value = 0.6449699222604717
$ws.Cells.Item(3, 11).Value = 0.3699768206640927
$ws.Cells.Item(4, 2).Value = 0.05018831929411549
$ws.Cells.Item(4, 3).Value = 0.8037777881109285
$ws.Cells.Item(4, 4).Value = 0.8147464657516995
$ws.Cells.Item(4, 5).Value = 0.4345783383908784
$ws.Cells.Item(4, 6).Value = 0.3855269040122699
$ws.Cells.Item(4, 7).Value = 0.6441032906516032
$ws.Cells.Item(4, 8).Value = 0.4725811863294772
$ws.Cells.Item(4, 9).Value = 0.631351887111854
$ws.Cells.Item(4, 10).Value = 0.3514560234789169
$ws.Cells.Item(4, 11).Value = 0.1878533956106533
$ws.Cells.Item(5, 2).Value = 0.761457155790539
$ws.Cells.Item(5, 3).Value = 0.7855419918959663
$ws.Cells.Item(5, 4).Value = 0.4315577385859573
$ws.Cells.Item(5, 5).Value = 0.3666656204007895
$ws.Cells.Item(5, 6).Value = 0.623772159586399
$ws.Cells.Item(5, 7).Value = 0.4580510308952294
$ws.Cells.Item(5, 8).Value = 0.6147108416195859
$ws.Cells.Item(5, 9).Value = 0.3338721780501667
$ws.Cells.Item(5, 10).Value = 0.1713795120806105
$ws.Cells.Item(5, 11).Value = 0.04814754147425859
$ws.Cells.Item(6, 2).Value = 1.125739732513876
$ws.Cells.Item(6, 3).Value = 0.5069288003698943
$ws.Cells.Item(6, 4).Value = 0.1755723340877904
$ws.Cells.Item(6, 5).Value = 0.6484660721731321
$ws.Cells.Item(6, 6).Value = 0.4674304344515254
$ws.Cells.Item(6, 7).Value = 0.5599000430081948
$ws.Cells.Item(6, 8).Value = 0.3149780871614969
$ws.Cells.Item(6, 9).Value = 0.1565155112655648
$ws.Cells.Item(6, 10).Value = 0.01943583084488787
$ws.Cells.Item(6, 11).Value = 0.5403464745801891
$ws.Cells.Item(7, 2).Value = 0.9576236723601449
$ws.Cells.Item(7, 3).Value = 0.2226830363001488
$ws.Cells.Item(7, 4).Value = 0.4081828413823239
$ws.Cells.Item(7, 5).Value = 0.5033755155173882
$ws.Cells.Item(7, 6).Value = 0.557849285568602
$ws.Cells.Item(7, 7).Value = 0.2384098477254667
$ws.Cells.Item(7, 8).Value = 0.1291161168902814
$ws.Cells.Item(7, 9).Value = -0.006853219352837503
$ws.Cells.Item(7, 10).Value = 0.4969491838668565
$ws.Cells.Item(7, 11).Value = 0.2970525035592049
$ws.Cells.Item(8, 2).Value = 0.5350086232236873
$ws.Cells.Item(8, 3).Value = 0.5414715811178482
$ws.Cells.Item(8, 4).Value = 0.3242970937011966
$ws.Cells.Item(8, 5).Value = 0.5860435243391693
$ws.Cells.Item(8, 6).Value = 0.2744351589178602
$ws.Cells.Item(8, 7).Value = 0.09157454668606776
$ws.Cells.Item(8, 8).Value = -0.01239555844019874
$ws.Cells.Item(8, 9).Value = 0.5011245128056051
$ws.Cells.Item(8, 10).Value = 0.2858677898194339
$ws.Cells.Item(9, 2).Value = 0.7770620722256893
$ws.Cells.Item(9, 3).Value = 0.4089758345755062
$ws.Cells.Item(9, 4).Value = 0.4404431577054521
$ws.Cells.Item(9, 5).Value = 0.2843889017771028
$ws.Cells.Item(9, 6).Value = 0.1069397163856127
$ws.Cells.Item(9, 7).Value = -0.05124389339762181
$ws.Cells.Item(9, 8).Value = 0.4852787037784192
$ws.Cells.Item(9, 9).Value = 0.2775335613519331
$ws.Cells.Item(10, 2).Value = 0.720000021056889
$ws.Cells.Item(10, 3).Value = 0.5575319321166473
$ws.Cells.Item(10, 4).Value = 0.121924760327536
$ws.Cells.Item(10, 5).Value = 0.1358986668508143
$ws.Cells.Item(10, 6).Value = -0.01550102094271821
$ws.Cells.Item(10, 7).Value = 0.4539510573947921
$ws.Cells.Item(10, 8).Value = 0.2743085116504074
$ws.Cells.Item(11, 2).Value = 0.8043571095207618
$ws.Cells.Item(11, 3).Value = 0.1395806278654255
$ws.Cells.Item(11, 4).Value = 0.04127642773981979
$ws.Cells.Item(11, 5).Value = 0.01671636730576187
$ws.Cells.Item(11, 6).Value = 0.4663391832225094
$ws.Cells.Item(11, 7).Value = 0.2534447081011285
$ws.Cells.Item(12, 2).Value = 0.3792569872816374
$ws.Cells.Item(12, 3).Value = 0.1263006712090333
$ws.Cells.Item(12, 4).Value = -0.09976821105225223
$ws.Cells.Item(12, 5).Value = 0.4814444548743619
$ws.Cells.Item(12, 6).Value = 0.2766837437271186
$ws.Cells.Item(13, 2).Value = 0.2911607482431205
$ws.Cells.Item(13, 3).Value = -0.08623093880433602
$ws.Cells.Item(13, 4).Value = 0.4184715358843989
$ws.Cells.Item(13, 5).Value = 0.2867219094086165
$ws.Cells.Item(14, 2).Value = 0.167618103061072
$ws.Cells.Item(14, 3).Value = 0.5177895860664353
$ws.Cells.Item(14, 4).Value = 0.1751453671933744
$ws.Cells.Item(15, 2).Value = 0.5618492773058843
$ws.Cells.Item(15, 3).Value = 0.1965658720679752
$ws.Cells.Item(16, 2).Value = 0.4328090033804217
